$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Battery_Data")
$ws2 = $wb.Worksheets.Item("Yearly BRC")

# --- Sheet1 (Battery_Data): update B2:B5 values ---
$ws1.Range("B2").Value = 2712216.28371
$ws1.Range("B3").Value = 2007040.0499454
$ws1.Range("B4").Value = 40140.800998908
$ws1.Range("B5").Value = 664663.574098

# --- Sheet2 (Yearly BRC): update existing B2:B3 values ---
$ws2.Range("B2").Value = 47443.72683102623
$ws2.Range("B3").Value = 47632.94166258632

# --- Sheet2: append new rows 4-21 (Battery Replacement Cost at y = 3..20) ---
$rows = @()

$rows += ,@(4, "Battery Replacement Cost at y = 3", 47618.10862433855)
$rows += ,@(5, "Battery Replacement Cost at y = 4", 47619.97206547673)
$rows += ,@(6, "Battery Replacement Cost at y = 5", 47624.80395492163)
$rows += ,@(7, "Battery Replacement Cost at y = 6", 47627.06485060284)
$rows += ,@(8, "Battery Replacement Cost at y = 7", 47628.36621991211)
$rows += ,@(9, "Battery Replacement Cost at y = 8", 47629.33247145954)
$rows += ,@(10, "Battery Replacement Cost at y = 9", 47630.27460722077)
$rows += ,@(11, "Battery Replacement Cost at y = 10", 47631.21703531164)
$rows += ,@(12, "Battery Replacement Cost at y = 11", 47632.18524322069)
$rows += ,@(13, "Battery Replacement Cost at y = 12", 47633.22584214556)
$rows += ,@(14, "Battery Replacement Cost at y = 13", 47634.17763780378)
$rows += ,@(15, "Battery Replacement Cost at y = 14", 47635.27687057104)
$rows += ,@(16, "Battery Replacement Cost at y = 15", 47636.40274215503)
$rows += ,@(17, "Battery Replacement Cost at y = 16", 47637.68278665161)
$rows += ,@(18, "Battery Replacement Cost at y = 17", 47638.98149000909)
$rows += ,@(19, "Battery Replacement Cost at y = 18", 47640.27735621159)
$rows += ,@(20, "Battery Replacement Cost at y = 19", 47641.43246232736)
$rows += ,@(21, "Battery Replacement Cost at y = 20", 47600.21489230519)

foreach ($row in $rows) {
    $r = $row[0]
    $label = $row[1]
    $val = $row[2]

    $srcA = $ws2.Range("A3")
    $dstA = $ws2.Range("A$r")
    $srcA.Copy($dstA)
    $dstA.Value = $label

    $ws2.Range("B$r").Value = $val
}

# --- Active sheet / tab selection: Battery_Data becomes the active tab ---
$ws1.Activate()

Write-Host "Edit complete"

